# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows for "Femacal de La Calera" / Frutilla
# right after the existing row 495, pushing the remaining data (old rows
# 496-535) down to rows 498-537 (matches dimension A1:T535 -> A1:T537).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 496-497 (existing rows shift down by 2).
$ws.Range("A496:A497").EntireRow.Insert()

# --- New row 496 ---
$ws.Range("A496").Value = 3
$ws.Range("B496").Value = "Femacal de La Calera"
$ws.Range("C496").Value = "Coquimbo"
$ws.Range("D496").Value = 45021
$ws.Range("E496").Value = 5
$ws.Range("F496").Value = "Fruta"
$ws.Range("G496").Value = 100101
$ws.Range("H496").Value = "Berries"
$ws.Range("I496").Value = 100112025
$ws.Range("J496").Value = "Frutilla"
$ws.Range("K496").Value = "Sin especificar"
$ws.Range("L496").Value = "Primera"
$ws.Range("M496").Value = 50
$ws.Range("N496").Value = 7000
$ws.Range("O496").Value = 7000
$ws.Range("P496").Value = 7000
$ws.Range("Q496").Value = "$/bandeja 7 kilos"
$ws.Range("R496").Value = "Provincia de Melipilla"
$ws.Range("S496").Value = 1000
$ws.Range("T496").Value = 7

# --- New row 497 ---
$ws.Range("A497").Value = 3
$ws.Range("B497").Value = "Femacal de La Calera"
$ws.Range("C497").Value = "Coquimbo"
$ws.Range("D497").Value = 45021
$ws.Range("E497").Value = 5
$ws.Range("F497").Value = "Fruta"
$ws.Range("G497").Value = 100101
$ws.Range("H497").Value = "Berries"
$ws.Range("I497").Value = 100112025
$ws.Range("J497").Value = "Frutilla"
$ws.Range("K497").Value = "Sin especificar"
$ws.Range("L497").Value = "Segunda"
$ws.Range("M497").Value = 40
$ws.Range("N497").Value = 6000
$ws.Range("O497").Value = 6000
$ws.Range("P497").Value = 6000
$ws.Range("Q497").Value = "$/bandeja 7 kilos"
$ws.Range("R497").Value = "Provincia de Melipilla"
$ws.Range("S497").Value = 857
$ws.Range("T497").Value = 7
